$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 - this shifts the existing rows 4-9 down to 5-10,
# carrying their formatting (incl. the date number format on column D) with them.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with a fresh weekly record (same market /
# category / etc. as its neighbours, new date 2021-11-10 = serial 44512).
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44512
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 300000000
$ws.Range("G4").Value = "Espárragos"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 950
$ws.Range("N4").Value = "$/kilo"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 950
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"
